# This script reproduces the "minor changes to Name class and subsequent
# functions" commit: every single-keyword cell in columns H (keyword1) and
# J (keyword2) is rewritten from a bare word (e.g. "Hire") into a small
# Python-list-repr string that also records the keyword's part of speech
# (e.g. "['Hire','verb']"), matching the part of speech implied by column
# B's "x + y" pattern label for that row.
#
# Values are written in the exact order the new shared-strings are expected
# to be appended to the workbook's shared string table (verb-tagged pairs
# first in row order, then noun-tagged pairs, then the single suffix pair),
# so the regenerated xl/sharedStrings.xml lines up with the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# verb + verb / verb + noun rows: keyword1/keyword2 tagged 'verb'
$ws.Cells.Item(2, 8).Value  = "['Hire','verb']"
$ws.Cells.Item(2, 10).Value = "['Failed','verb']"
$ws.Cells.Item(4, 8).Value  = "['Offers','verb']"
$ws.Cells.Item(4, 10).Value = "['Vary','verb']"
$ws.Cells.Item(5, 8).Value  = "['Flowing','verb']"
$ws.Cells.Item(5, 10).Value = "['Guide','verb']"
$ws.Cells.Item(7, 8).Value  = "['Build','verb']"
$ws.Cells.Item(9, 10).Value = "['Flame','verb']"
$ws.Cells.Item(10, 8).Value = "['Output','verb']"
$ws.Cells.Item(11, 8).Value = "['Mind','verb']"

# noun + and/to + noun rows, and the noun half of verb + noun rows
$ws.Cells.Item(3, 8).Value   = "['Media','noun']"
$ws.Cells.Item(3, 10).Value  = "['Then','noun']"
$ws.Cells.Item(6, 8).Value   = "['Women','noun']"
$ws.Cells.Item(6, 10).Value  = "['Gaea','noun']"
$ws.Cells.Item(7, 10).Value  = "['Lessons','noun']"
$ws.Cells.Item(8, 8).Value   = "['Area','noun']"
$ws.Cells.Item(10, 10).Value = "['Family','noun']"
$ws.Cells.Item(11, 10).Value = "['Adobe','noun']"

# noun + suffix row
$ws.Cells.Item(8, 10).Value = "['onym','suffix']"

# Same word/tag pair recurring later in the sheet (row 9, col H) re-uses the
# already-interned shared string instead of creating a duplicate.
$ws.Cells.Item(9, 8).Value = "['Vary','verb']"

# Widen a handful of columns so the long list-repr strings stay legible.
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666   # -> ~15
$ws.Columns.Item(3).ColumnWidth = 8.0                  # -> ~8.83 (default-ish)
$ws.Columns.Item(4).ColumnWidth = 8.0
$ws.Columns.Item(5).ColumnWidth = 8.0
$ws.Columns.Item(6).ColumnWidth = 8.0
$ws.Columns.Item(7).ColumnWidth = 8.0
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666   # -> ~15.5
$ws.Columns.Item(9).ColumnWidth = 8.0
$ws.Columns.Item(10).ColumnWidth = 15.83               # -> ~16.66

# Move the active selection.
$ws.Range("J22").Select()
